# Apply the "Add files via upload" change:
# - drop the old "tag6" column (H) and insert a new "CM"-tag column at C
# - retag row 1 (header) label from "ling" to "link", add "CM" header at C1
# - add a "C"/"M" tag value in the new column for the existing company row
#   and for a newly added sibling company row (Social/Default.aspx)
# - move the cell selection to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old tag6 column, then insert a fresh blank column at C.
# This shifts C..G -> D..H while leaving I/J (interlinkregex/finallinkregex) in place,
# matching the diff (dimension stays A1:J6).
$ws.Columns("H:H").Delete()
$ws.Columns("C:C").Insert()

# --- Row 1 (header row) ---
$ws.Range("A1").Value2 = "link"
$ws.Range("C1").Value2 = "CM"

# --- Row 2 (existing company row) ---
$ws.Range("C2").Value2 = "C"

# --- Row 3 (new company row, mirrors row 2 but with a different link & tag) ---
$ws.Range("A3").Value2 = "http://career.cmbchina.com/Social/Default.aspx"
$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = "M"
$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2
$ws.Range("I3").Value2 = $ws.Range("I2").Value2
$ws.Range("J3").Value2 = $ws.Range("J2").Value2

# --- Move active selection to D13, as recorded in the saved view state ---
$ws.Range("D13").Select()
